$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 and add new rows 8-10 with refreshed TPM data
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Qrfp"
$ws.Range("C2").Value = "P2ry14"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2329866666666667
$ws.Range("H2").Value = 0.69896
$ws.Range("I2").Value = 0.4906425117280099
$ws.Range("J2").Value = 0.4906425117280099
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.047064
$ws.Range("N2").Value = 0.141192
$ws.Range("O2").Value = 0.003291309722706829
$ws.Range("P2").Value = 0.00329130972270683
$ws.Range("Q2").Value = 0.01096528448
$ws.Range("R2").Value = 0.09868756032000001
$ws.Range("S2").Value = 0.001614856469223699
$ws.Range("T2").Value = 0.001614856469223699

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Qrfp"
$ws.Range("C3").Value = "P2ry14"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2329866666666667
$ws.Range("H3").Value = 0.69896
$ws.Range("I3").Value = 0.4906425117280099
$ws.Range("J3").Value = 0.4906425117280099
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.005187333333334
$ws.Range("N3").Value = 27.015562
$ws.Range("O3").Value = 0.629756515064516
$ws.Range("P3").Value = 0.6297565150645161
$ws.Range("Q3").Value = 2.098088579502222
$ws.Range("R3").Value = 18.88279721552
$ws.Range("S3").Value = 0.3089853183283324
$ws.Range("T3").Value = 0.3089853183283325

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Qrfp"
$ws.Range("C4").Value = "P2ry14"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2329866666666667
$ws.Range("H4").Value = 0.69896
$ws.Range("I4").Value = 0.4906425117280099
$ws.Range("J4").Value = 0.4906425117280099
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.247223333333333
$ws.Range("N4").Value = 15.74167
$ws.Range("O4").Value = 0.366952175212777
$ws.Range("P4").Value = 0.3669521752127771
$ws.Range("Q4").Value = 1.222533073688889
$ws.Range("R4").Value = 11.0027976632
$ws.Range("S4").Value = 0.1800423369304537
$ws.Range("T4").Value = 0.1800423369304537

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Qrfp"
$ws.Range("C5").Value = "P2ry14"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09243566666666668
$ws.Range("H5").Value = 0.277307
$ws.Range("I5").Value = 0.1946586399790535
$ws.Range("J5").Value = 0.1946586399790535
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.047064
$ws.Range("N5").Value = 0.141192
$ws.Range("O5").Value = 0.003291309722706829
$ws.Range("P5").Value = 0.00329130972270683
$ws.Range("Q5").Value = 0.004350392216000001
$ws.Range("R5").Value = 0.03915352994400001
$ws.Range("S5").Value = 0.0006406818743719472
$ws.Range("T5").Value = 0.0006406818743719473

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Qrfp"
$ws.Range("C6").Value = "P2ry14"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.09243566666666668
$ws.Range("H6").Value = 0.277307
$ws.Range("I6").Value = 0.1946586399790535
$ws.Range("J6").Value = 0.1946586399790535
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.005187333333334
$ws.Range("N6").Value = 27.015562
$ws.Range("O6").Value = 0.629756515064516
$ws.Range("P6").Value = 0.6297565150645161
$ws.Range("Q6").Value = 0.8324004946148891
$ws.Range("R6").Value = 7.491604451534
$ws.Range("S6").Value = 0.122587546740407
$ws.Range("T6").Value = 0.122587546740407

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Qrfp"
$ws.Range("C7").Value = "P2ry14"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.09243566666666668
$ws.Range("H7").Value = 0.277307
$ws.Range("I7").Value = 0.1946586399790535
$ws.Range("J7").Value = 0.1946586399790535
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.247223333333333
$ws.Range("N7").Value = 15.74167
$ws.Range("O7").Value = 0.366952175212777
$ws.Range("P7").Value = 0.3669521752127771
$ws.Range("Q7").Value = 0.4850305869655556
$ws.Range("R7").Value = 4.36527528269
$ws.Range("S7").Value = 0.07143041136427453
$ws.Range("T7").Value = 0.07143041136427454

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Qrfp"
$ws.Range("C8").Value = "P2ry14"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.149438
$ws.Range("H8").Value = 0.448314
$ws.Range("I8").Value = 0.3146988482929367
$ws.Range("J8").Value = 0.3146988482929367
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.047064
$ws.Range("N8").Value = 0.141192
$ws.Range("O8").Value = 0.003291309722706829
$ws.Range("P8").Value = 0.00329130972270683
$ws.Range("Q8").Value = 0.007033150032
$ws.Range("R8").Value = 0.063298350288
$ws.Range("S8").Value = 0.001035771379111184
$ws.Range("T8").Value = 0.001035771379111184

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Qrfp"
$ws.Range("C9").Value = "P2ry14"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.149438
$ws.Range("H9").Value = 0.448314
$ws.Range("I9").Value = 0.3146988482929367
$ws.Range("J9").Value = 0.3146988482929367
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.005187333333334
$ws.Range("N9").Value = 27.015562
$ws.Range("O9").Value = 0.629756515064516
$ws.Range("P9").Value = 0.6297565150645161
$ws.Range("Q9").Value = 1.345717184718667
$ws.Range("R9").Value = 12.111454662468
$ws.Range("S9").Value = 0.1981836499957766
$ws.Range("T9").Value = 0.1981836499957767

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Qrfp"
$ws.Range("C10").Value = "P2ry14"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.149438
$ws.Range("H10").Value = 0.448314
$ws.Range("I10").Value = 0.3146988482929367
$ws.Range("J10").Value = 0.3146988482929367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.247223333333333
$ws.Range("N10").Value = 15.74167
$ws.Range("O10").Value = 0.366952175212777
$ws.Range("P10").Value = 0.3669521752127771
$ws.Range("Q10").Value = 0.7841345604866666
$ws.Range("R10").Value = 7.05721104438
$ws.Range("S10").Value = 0.1154794269180488
$ws.Range("T10").Value = 0.1154794269180489

